$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases Import")

# Only the header row is restructured: a new "Scenario" column is inserted
# after Description, pushing Step Number/Action/Expected Result one column
# to the right (E->F, F->G, G->H). The existing data rows (2-9) keep their
# original A:G layout untouched - write the shifted header cells directly
# instead of a real column insert (which would also push the data rows).
$g1 = $ws.Range("G1").Value2
$f1 = $ws.Range("F1").Value2
$e1 = $ws.Range("E1").Value2

$ws.Range("H1").Value = $g1
$ws.Range("G1").Value = $f1
$ws.Range("F1").Value = $e1
$ws.Range("E1").Value = "Scenario"

# Suite Name -> Submodule Name
$ws.Range("A1").Value = "Submodule Name"

# The header row loses its bold/blue "header" style in the target workbook -
# all header cells fall back to the default (unstyled) format.
$ws.Range("A1:H1").Style = "Normal"
